$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value2 = "53.603.80"
$ws.Cells.Item(2, 5).Value2 = "  -5.13%  "

# Row 3
$ws.Cells.Item(3, 4).Value2 = "2.211.01"
$ws.Cells.Item(3, 5).Value2 = "  -6.99%  "

# Row 4
$ws.Cells.Item(4, 4).Value2 = "'1.00"
$ws.Cells.Item(4, 5).Value2 = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).Value2 = "'486.72"
$ws.Cells.Item(5, 5).Value2 = "  -3.76%  "

# Row 6
$ws.Cells.Item(6, 4).Value2 = "'125.17"
$ws.Cells.Item(6, 5).Value2 = "  -4.24%  "

# Row 7
$ws.Cells.Item(7, 4).Value2 = "'0.995"
$ws.Cells.Item(7, 5).Value2 = "  -0.25%  "

# Row 8
$ws.Cells.Item(8, 5).Value2 = "  -4.32%  "

# Row 9
$ws.Cells.Item(9, 4).Value2 = "2.239.72"
$ws.Cells.Item(9, 5).Value2 = "  -6.27%  "

# Row 10
$ws.Cells.Item(10, 4).Value2 = "'0.0920"
$ws.Cells.Item(10, 5).Value2 = "  -6.73%  "

# Row 11
$ws.Cells.Item(11, 5).Value2 = "  -0.22%  "

# Row 12
$ws.Cells.Item(12, 4).Value2 = "'0.320"
$ws.Cells.Item(12, 5).Value2 = "  -2.62%  "

# Row 13
$ws.Cells.Item(13, 4).Value2 = "'4.62"
$ws.Cells.Item(13, 5).Value2 = "  -4.71%  "

# Row 14
$ws.Cells.Item(14, 4).Value2 = "2.610.13"
$ws.Cells.Item(14, 5).Value2 = "  -6.78%  "

# Row 15
$ws.Cells.Item(15, 5).Value2 = "  -1.91%  "

# Row 16
$ws.Cells.Item(16, 4).Value2 = "53.537.91"
$ws.Cells.Item(16, 5).Value2 = "  -5.17%  "

# Row 17
$ws.Cells.Item(17, 5).Value2 = "  -4.16%  "

# Row 18
$ws.Cells.Item(18, 4).Value2 = "2.234.32"
$ws.Cells.Item(18, 5).Value2 = "  -6.32%  "

# Row 19
$ws.Cells.Item(19, 5).Value2 = "  -4.14%  "

# Row 20
$ws.Cells.Item(20, 5).Value2 = "  -1.99%  "

# Row 21
$ws.Cells.Item(21, 4).Value2 = "'295.97"
$ws.Cells.Item(21, 5).Value2 = "  -4.42%  "

# Row 22
$ws.Cells.Item(22, 4).Value2 = "'6.17"
$ws.Cells.Item(22, 5).Value2 = "  -2.08%  "

# Row 23
$ws.Cells.Item(23, 4).Value2 = "'0.996"
$ws.Cells.Item(23, 5).Value2 = "  -0.39%  "

# Row 24
$ws.Cells.Item(24, 4).Value2 = "'63.64"
$ws.Cells.Item(24, 5).Value2 = "  -4.09%  "

# Row 25
$ws.Cells.Item(25, 4).Value2 = "'0.998"
$ws.Cells.Item(25, 5).Value2 = "  +0.13%  "

# Row 26
$ws.Cells.Item(26, 4).Value2 = "'0.367"
$ws.Cells.Item(26, 5).Value2 = "  -1.21%  "

# Row 27
$ws.Cells.Item(27, 2).Value2 = "WrappedeETH"
$ws.Cells.Item(27, 3).Value2 = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(27, 4).Value2 = "2.339.45"
$ws.Cells.Item(27, 5).Value2 = "  -6.14%  "

# Row 28
$ws.Cells.Item(28, 2).Value2 = "Kaspa"
$ws.Cells.Item(28, 3).Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(28, 4).Value2 = "'0.147"
$ws.Cells.Item(28, 5).Value2 = "  -0.92%  "

# Row 29
$ws.Cells.Item(29, 5).Value2 = "  -3.14%  "

# Row 30
$ws.Cells.Item(30, 4).Value2 = "'162.93"
$ws.Cells.Item(30, 5).Value2 = "  -6.20%  "

# Row 31
$ws.Cells.Item(31, 5).Value2 = "  -4.02%  "

# Row 33
$ws.Cells.Item(33, 5).Value2 = "  -6.27%  "

# Row 34
$ws.Cells.Item(34, 5).Value2 = "  -1.07%  "

# Row 35
$ws.Cells.Item(35, 4).Value2 = "'0.993"
$ws.Cells.Item(35, 5).Value2 = "  -0.30%  "

# Row 36
$ws.Cells.Item(36, 5).Value2 = "  -2.17%  "

# Row 37
$ws.Cells.Item(37, 4).Value2 = "'17.29"
$ws.Cells.Item(37, 5).Value2 = "  -2.28%  "

# Row 38
$ws.Cells.Item(38, 5).Value2 = "  -0.85%  "

# Row 39
$ws.Cells.Item(39, 4).Value2 = "'0.834"
$ws.Cells.Item(39, 5).Value2 = "  +1.30%  "

# Row 40
$ws.Cells.Item(40, 5).Value2 = "  -3.78%  "

# Row 41
$ws.Cells.Item(41, 4).Value2 = "'35.18"
$ws.Cells.Item(41, 5).Value2 = "  -3.64%  "

# Row 42
$ws.Cells.Item(42, 5).Value2 = "  -0.91%  "

# Row 43
$ws.Cells.Item(43, 5).Value2 = "  -1.41%  "

# Row 44
$ws.Cells.Item(44, 2).Value2 = "Aave"
$ws.Cells.Item(44, 3).Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(44, 4).Value2 = "'127.68"
$ws.Cells.Item(44, 5).Value2 = "  -0.52%  "

# Row 45
$ws.Cells.Item(45, 2).Value2 = "Filecoin"
$ws.Cells.Item(45, 3).Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45, 4).Value2 = "'3.30"
$ws.Cells.Item(45, 5).Value2 = "  -2.68%  "

# Row 46
$ws.Cells.Item(46, 5).Value2 = "  -2.77%  "

# Row 47
$ws.Cells.Item(47, 5).Value2 = "  -1.79%  "

# Row 48
$ws.Cells.Item(48, 4).Value2 = "'0.536"
$ws.Cells.Item(48, 5).Value2 = "  -5.29%  "

# Row 49
$ws.Cells.Item(49, 4).Value2 = "'234.84"
$ws.Cells.Item(49, 5).Value2 = "  -2.39%  "

# Row 50
$ws.Cells.Item(50, 4).Value2 = "'0.0472"
$ws.Cells.Item(50, 5).Value2 = "  -2.10%  "

# Row 51
$ws.Cells.Item(51, 5).Value2 = "  -3.41%  "
